# Autosize nodes to image:
# Shift every node's "y" coordinate (column C, rows 2-45) up by 1,
# and fix the x coordinate (column B) for node ids 22 and 23 (rows 23-24).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("nodes")

# Fix x (column B) values for the two nodes whose deployment slot moved.
$ws.Cells.Item(23, 2).Value = 4
$ws.Cells.Item(24, 2).Value = 6

# Shift y (column C) values for every data row by +1.
for ($row = 2; $row -le 45; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value()
    $cell.Value = $current + 1
}

# Update the view state to match the saved selection/scroll position.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("J23").Select()
